$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2764.0908
$ws.Range("I94").Value = 2667.2222
$ws.Range("K94").Value = 2667.2222
$ws.Range("M94").Value = -2216.2222

$ws.Range("H107").Value = 443.53845
$ws.Range("I107").Value = 443.53845
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 443.53845
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1476.46155
$ws.Range("N107").ClearContents()

$ws.Range("H132").Value = 5495759
$ws.Range("I132").Value = 6212297
$ws.Range("K132").Value = 18636891
$ws.Range("M132").Value = -18634361

$ws.Range("H135").Value = 1232.1538
$ws.Range("I135").Value = 1177.9
$ws.Range("J135").Value = 1413
$ws.Range("K135").Value = 10601.1
$ws.Range("L135").Value = 12717
$ws.Range("M135").Value = -8066.1
$ws.Range("N135").Value = -17787

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1201.5454
$ws.Range("I2").Value = 901.75
$ws.Range("J2").Value = 1561.3
$ws.Range("K2").Value = 901.75
$ws.Range("L2").Value = 1561.3
$ws.Range("M2").Value = -788.75
$ws.Range("N2").Value = -1787.3

$ws.Range("H22").Value = 4328.5713
$ws.Range("I22").Value = 1060
$ws.Range("K22").Value = 1060
$ws.Range("M22").Value = -761

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H74").Value = 1726.579
$ws.Range("I74").Value = 1665
$ws.Range("J74").Value = 2250
$ws.Range("K74").Value = 1665
$ws.Range("L74").Value = 2250
$ws.Range("M74").Value = -791
$ws.Range("N74").Value = -3998

$ws.Range("H77").Value = 1726.579
$ws.Range("I77").Value = 1665
$ws.Range("J77").Value = 2250
$ws.Range("K77").Value = 8325
$ws.Range("L77").Value = 11250
$ws.Range("M77").Value = -3957
$ws.Range("N77").Value = -19986

$ws.Range("H116").Value = 1201.5454
$ws.Range("I116").Value = 901.75
$ws.Range("J116").Value = 1561.3
$ws.Range("K116").Value = 901.75
$ws.Range("L116").Value = 1561.3
$ws.Range("M116").Value = 1392.25
$ws.Range("N116").Value = -6149.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1201.5454
$ws.Range("I3").Value = 901.75
$ws.Range("J3").Value = 1561.3
$ws.Range("K3").Value = 901.75
$ws.Range("L3").Value = 1561.3
$ws.Range("M3").Value = -787.75
$ws.Range("N3").Value = -1789.3

$ws.Range("H88").Value = 29089.092
$ws.Range("J88").Value = 29089.092
$ws.Range("L88").Value = 29089.092
$ws.Range("N88").Value = -29901.092

$ws.Range("H91").Value = 29089.092
$ws.Range("J91").Value = 29089.092
$ws.Range("L91").Value = 29089.092
$ws.Range("N91").Value = -31897.092

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1876
$ws.Range("I16").Value = 1133.3334
$ws.Range("J16").Value = 2990
$ws.Range("K16").Value = 1133.3334
$ws.Range("L16").Value = 2990
$ws.Range("M16").Value = -846.3334
$ws.Range("N16").Value = -3564

$ws.Range("H31").Value = 2599144.2
$ws.Range("I31").Value = 1656
$ws.Range("J31").Value = 8001920
$ws.Range("K31").Value = 1656
$ws.Range("L31").Value = 8001920
$ws.Range("M31").Value = -1361
$ws.Range("N31").Value = -8002510

$ws.Range("H34").Value = 2599144.2
$ws.Range("I34").Value = 1656
$ws.Range("J34").Value = 8001920
$ws.Range("K34").Value = 1656
$ws.Range("L34").Value = 8001920
$ws.Range("M34").Value = -1454
$ws.Range("N34").Value = -8002324

$ws.Range("H107").Value = 539.1875
$ws.Range("I107").Value = 539.1875
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 539.1875
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1380.8125
$ws.Range("N107").ClearContents()

$ws.Range("H113").Value = 1876
$ws.Range("I113").Value = 1133.3334
$ws.Range("J113").Value = 2990
$ws.Range("K113").Value = 1133.3334
$ws.Range("L113").Value = 2990
$ws.Range("M113").Value = 1036.6666
$ws.Range("N113").Value = -7330

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 626.9167
$ws.Range("I5").Value = 442.6842
$ws.Range("J5").Value = 1327
$ws.Range("K5").Value = 1328.0526
$ws.Range("L5").Value = 3981
$ws.Range("M5").Value = -1216.0526
$ws.Range("N5").Value = -4205

$ws.Range("H75").Value = 1403
$ws.Range("I75").Value = 1537.3334
$ws.Range("J75").Value = 1000
$ws.Range("K75").Value = 4612.0002
$ws.Range("L75").Value = 3000
$ws.Range("M75").Value = -3614.0002
$ws.Range("N75").Value = -4996

$ws.Range("H78").Value = 1403
$ws.Range("I78").Value = 1537.3334
$ws.Range("J78").Value = 1000
$ws.Range("K78").Value = 13836.0006
$ws.Range("L78").Value = 9000
$ws.Range("M78").Value = -8844.000599999999
$ws.Range("N78").Value = -18984

$ws.Range("H99").Value = 10952.083
$ws.Range("I99").Value = 925
$ws.Range("J99").Value = 11863.637
$ws.Range("K99").Value = 2775
$ws.Range("L99").Value = 35590.911
$ws.Range("M99").Value = -529
$ws.Range("N99").Value = -40082.911

$ws.Range("H122").Value = 2150.1
$ws.Range("I122").Value = 2920.6
$ws.Range("K122").Value = 26285.4
$ws.Range("M122").Value = -23835.4

$ws.Range("H123").Value = 2320
$ws.Range("I123").Value = 1866.6666
$ws.Range("K123").Value = 5599.9998
$ws.Range("M123").Value = -3149.9998

$ws.Range("H124").Value = 5000
$ws.Range("I124").Value = 1000
$ws.Range("K124").Value = 3000
$ws.Range("M124").Value = 1910

$ws.Range("H125").Value = 2643.3333
$ws.Range("J125").Value = 6900
$ws.Range("L125").Value = 20700
$ws.Range("N125").Value = -30540

$ws.Range("H131").Value = 3741180.2
$ws.Range("J131").Value = 10192614
$ws.Range("L131").Value = 30577842
$ws.Range("N131").Value = -30587922

$ws.Range("H135").Value = 626.9167
$ws.Range("I135").Value = 442.6842
$ws.Range("J135").Value = 1327
$ws.Range("K135").Value = 3984.1578
$ws.Range("L135").Value = 11943
$ws.Range("M135").Value = -1449.1578
$ws.Range("N135").Value = -17013

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 625.8333
$ws.Range("I22").Value = 264.2857
$ws.Range("J22").Value = 1132
$ws.Range("K22").Value = 264.2857
$ws.Range("L22").Value = 1132
$ws.Range("M22").Value = 30.71429999999998
$ws.Range("N22").Value = -1722

$ws.Range("H27").Value = 625.8333
$ws.Range("I27").Value = 264.2857
$ws.Range("J27").Value = 1132
$ws.Range("K27").Value = 264.2857
$ws.Range("L27").Value = 1132
$ws.Range("M27").Value = -157.2857
$ws.Range("N27").Value = -1346

$ws.Range("H46").Value = 1715.9032
$ws.Range("I46").Value = 1103.9
$ws.Range("J46").Value = 2007.3334
$ws.Range("K46").Value = 1103.9
$ws.Range("L46").Value = 2007.3334
$ws.Range("M46").Value = -915.9000000000001
$ws.Range("N46").Value = -2383.3334
